# "added screenshot in extent reports"
# The underlying data change: the old "jR@1372jitu" credential value is
# replaced by "N" (SignIn/SignUp "runmode" columns) while the
# SearchProductTest sheet's hyperlinked "jR@1372jitu" cells are replaced
# with two new placeholder values ("hhfgfgf" / "kjjkhjghfgf"); the mailto
# hyperlinks on those cells are kept, with their original display text
# cached on the hyperlink itself.

$wb = $excel.ActiveWorkbook

$wsSignIn  = $wb.Worksheets.Item("SignInTest")
$wsSignUp  = $wb.Worksheets.Item("SignUpTest")
$wsSearch  = $wb.Worksheets.Item("SearchProductTest")

# --- Update cell values -----------------------------------------------
# Order matters: it determines how new entries land in the shared-string
# table, matching the order they were typed in the original edit.
$wsSignIn.Range("C3").Value = "N"
$wsSignUp.Range("B2").Value = "N"
$wsSearch.Range("E3").Value = "N"
$wsSearch.Range("D3").Value = "hhfgfgf"
$wsSearch.Range("D2").Value = "kjjkhjghfgf"

# Re-point the cached hyperlink display text back to the original
# "jR@1372jitu" value (Excel keeps this around even after the visible
# cell text is changed, since the hyperlink target itself is untouched).
foreach ($hl in $wsSearch.Hyperlinks) {
    $hl.TextToDisplay = "jR@1372jitu"
}

# --- Update the active-cell selections on each sheet -------------------
# Select on the non-final sheets first, then finish on SignInTest so it
# remains the active tab, matching the saved workbook view state.
$wsSignUp.Range("B2").Select()
$wsSearch.Range("E13").Select()
$wsSignIn.Range("C3").Select()
